$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "16/05/2021"
$ws.Range("C6").Value = "15/05/2024"
